$wb = $excel.ActiveWorkbook

# --- Existing "AddEmployee" sheet: selection moves from G2 to D6 ---
$ws1 = $wb.Worksheets.Item("AddEmployee")
[void]$ws1.Range("D6").Select()

# --- Add the new "AddCandidate" sheet after "AddEmployee" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "AddCandidate"

# Header row
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Email"
$ws2.Range("C1").Value = "Id"
$ws2.Range("D1").Value = "DOB"

# Row 2
$ws2.Range("A2").Value = "Auto Test"
$ws2.Range("B2").Value = "autotest"
$ws2.Range("C2").Value = "c@gmail.com"
$ws2.Range("D2").Value = 19091997

# Row 3
$ws2.Range("A3").Value = "AutoinvalidTest123"
$ws2.Range("C3").Value = "3@gmail"
$ws2.Range("D3").Value = 78
$ws2.Range("E3").Value = 882323098

# Hyperlinks on the email/id column, styled like the existing "Hyperlink" cells
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:c@gmail.com")
$ws2.Range("C2").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("C3"), "mailto:3@gmail")
$ws2.Range("C3").Style = "Hyperlink"

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 25.8
$ws2.Columns.Item(3).ColumnWidth = 11.17
$ws2.Columns.Item(5).ColumnWidth = 9

# Final selection on the new sheet + make it the active tab
[void]$ws2.Range("E3").Select()
[void]$ws2.Activate()
